$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws4 = $wb.Worksheets.Item(4)

# Sheet1 (展览)
$ws1.Range("F2").Value = 1765
$ws1.Range("F3").Value = 770
$ws1.Range("F4").Value = 506
$ws1.Range("F5").Value = 252
$ws1.Range("F6").Value = 458
$ws1.Range("F7").Value = 1091
$ws1.Range("F8").Value = 313
$ws1.Range("F9").Value = 17
$ws1.Range("F10").Value = 102
$ws1.Range("F12").Value = 1081
$ws1.Range("F14").Value = 90
$ws1.Range("F15").Value = 719
$ws1.Range("F16").Value = 780
$ws1.Range("F17").Value = 166
$ws1.Range("F18").Value = 22
$ws1.Range("F19").Value = 49
$ws1.Range("F20").Value = 612
$ws1.Range("F21").Value = 101
$ws1.Range("F22").Value = 1685
$ws1.Range("F23").Value = 1883
$ws1.Range("F24").Value = 475
$ws1.Range("F25").Value = 54
$ws1.Range("F26").Value = 1715
$ws1.Range("F27").Value = 254
$ws1.Range("F28").Value = 2522
$ws1.Range("F29").Value = 446
$ws1.Range("F31").Value = 648
$ws1.Range("F33").Value = 83
$ws1.Range("F35").Value = 866
$ws1.Range("F36").Value = 1576
$ws1.Range("F37").Value = 261
$ws1.Range("F39").Value = 507
$ws1.Range("F40").Value = 103
$ws1.Range("F41").Value = 95
$ws1.Range("F42").Value = 135

# Sheet2 (演出)
$ws2.Range("F4").Value = 118
$ws2.Range("F10").Value = 3
$ws2.Range("F12").Value = 59

# Sheet4 (全部类型)
$ws4.Range("F2").Value = 1765
$ws4.Range("F5").Value = 770
$ws4.Range("F6").Value = 506
$ws4.Range("F7").Value = 252
$ws4.Range("F8").Value = 458
$ws4.Range("F9").Value = 1091
$ws4.Range("F10").Value = 313
$ws4.Range("F11").Value = 17
$ws4.Range("F12").Value = 102
$ws4.Range("F14").Value = 1081
$ws4.Range("F16").Value = 719
$ws4.Range("F17").Value = 780
$ws4.Range("F18").Value = 166
$ws4.Range("F19").Value = 118
$ws4.Range("F20").Value = 118
$ws4.Range("F22").Value = 22
$ws4.Range("F24").Value = 49
$ws4.Range("F25").Value = 612
$ws4.Range("F26").Value = 101
$ws4.Range("F27").Value = 1685
$ws4.Range("F28").Value = 1883
$ws4.Range("F29").Value = 475
$ws4.Range("F30").Value = 54
$ws4.Range("F32").Value = 2522
$ws4.Range("F33").Value = 446
$ws4.Range("F37").Value = 59
$ws4.Range("F38").Value = 648
$ws4.Range("F40").Value = 83
$ws4.Range("F42").Value = 866
$ws4.Range("F43").Value = 1576
$ws4.Range("F45").Value = 261
$ws4.Range("F46").Value = 507
$ws4.Range("F47").Value = 103
$ws4.Range("F48").Value = 95
$ws4.Range("F49").Value = 135
